# Applies the scheduled-runner market-data refresh described in the commit.
# All target cells are static (non-formula) values, so this is a pure value write.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1253.9412
$ws.Range("I9").Value = 1303
$ws.Range("K9").Value = 1303
$ws.Range("M9").Value = -1134

$ws.Range("H137").Value = 4894.4
$ws.Range("J137").Value = 4825.6665
$ws.Range("L137").Value = 14476.9995
$ws.Range("N137").Value = -19576.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5249.6206
$ws.Range("I61").Value = 5301.1665
$ws.Range("J61").Value = 5002.2
$ws.Range("K61").Value = 5301.1665
$ws.Range("L61").Value = 5002.2
$ws.Range("M61").Value = -5089.1665
$ws.Range("N61").Value = -5426.2

$ws.Range("H122").Value = 3297.761
$ws.Range("I122").Value = 3258.6216
$ws.Range("J122").Value = 3458.6667
$ws.Range("K122").Value = 9775.8648
$ws.Range("L122").Value = 10376.0001
$ws.Range("M122").Value = -7325.864799999999
$ws.Range("N122").Value = -15276.0001

$ws.Range("H136").Value = 5249.6206
$ws.Range("I136").Value = 5301.1665
$ws.Range("J136").Value = 5002.2
$ws.Range("K136").Value = 15903.4995
$ws.Range("L136").Value = 15006.6
$ws.Range("M136").Value = -13353.4995
$ws.Range("N136").Value = -20106.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4770.971
$ws.Range("I31").Value = 971
$ws.Range("J31").Value = 5557.1724
$ws.Range("K31").Value = 971
$ws.Range("L31").Value = 5557.1724
$ws.Range("M31").Value = -676
$ws.Range("N31").Value = -6147.1724

$ws.Range("H34").Value = 4770.971
$ws.Range("I34").Value = 971
$ws.Range("J34").Value = 5557.1724
$ws.Range("K34").Value = 971
$ws.Range("L34").Value = 5557.1724
$ws.Range("M34").Value = -769
$ws.Range("N34").Value = -5961.1724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4498.4
$ws.Range("J39").Value = 4887.222
$ws.Range("L39").Value = 14661.666
$ws.Range("N39").Value = -15249.666

$ws.Range("H55").Value = 48398.89
$ws.Range("J55").Value = 16948.75
$ws.Range("L55").Value = 50846.25
$ws.Range("N55").Value = -51200.25

$ws.Range("H121").Value = 63159.8
$ws.Range("I121").Value = 1000
$ws.Range("K121").Value = 3000
$ws.Range("M121").Value = -1690

$ws.Range("H128").Value = 237999.33
$ws.Range("I128").Value = 237999.33
$ws.Range("K128").Value = 713997.99
$ws.Range("M128").Value = -709017.99

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1860.125
$ws.Range("I97").Value = 1497.0834
$ws.Range("K97").Value = 1497.0834
$ws.Range("M97").Value = -1001.0834

$ws.Range("H132").Value = 37040560
$ws.Range("I132").Value = 52634804
$ws.Range("K132").Value = 157904412
$ws.Range("M132").Value = -157901882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5264.643
$ws.Range("I22").Value = 3441.8696
$ws.Range("J22").Value = 13649.4
$ws.Range("K22").Value = 3441.8696
$ws.Range("L22").Value = 13649.4
$ws.Range("M22").Value = -3146.8696
$ws.Range("N22").Value = -14239.4

$ws.Range("H27").Value = 5264.643
$ws.Range("I27").Value = 3441.8696
$ws.Range("J27").Value = 13649.4
$ws.Range("K27").Value = 3441.8696
$ws.Range("L27").Value = 13649.4
$ws.Range("M27").Value = -3334.8696
$ws.Range("N27").Value = -13863.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2309.7354
$ws.Range("I122").Value = 2329.125
$ws.Range("K122").Value = 6987.375
$ws.Range("M122").Value = -4537.375

$ws.Range("H132").Value = 3628542.8
$ws.Range("I132").Value = 5053718.5
$ws.Range("K132").Value = 15161155.5
$ws.Range("M132").Value = -15158625.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 39999.5
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 39999.5
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 39999.5
$ws.Range("N124").Value = -49819.5

$ws.Range("H125").Value = 75000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 75000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 75000
$ws.Range("N125").Value = -84840

$ws.Range("H127").Value = 99500
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 99500
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 99500
$ws.Range("N127").Value = -109420

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0

$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0

$ws.Range("H131").Value = 63495
$ws.Range("I131").Value = 70000
$ws.Range("J131").Value = 61326.668
$ws.Range("K131").Value = 70000
$ws.Range("L131").Value = 61326.668
$ws.Range("M131").Value = -64960
$ws.Range("N131").Value = -71406.668

$ws.Range("H132").Value = 3839.2
$ws.Range("I132").Value = 2997
$ws.Range("J132").Value = 4049.75
$ws.Range("K132").Value = 8991
$ws.Range("L132").Value = 12149.25
$ws.Range("M132").Value = -6461
$ws.Range("N132").Value = -17209.25

$ws.Range("H133").Value = 100000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 100000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0

$ws.Range("H136").Value = 68220910
$ws.Range("I136").Value = 41734856
$ws.Range("J136").Value = 100004170
$ws.Range("K136").Value = 125204568
$ws.Range("L136").Value = 300012510
$ws.Range("M136").Value = -125202018
$ws.Range("N136").Value = -300017610

$ws.Range("H137").Value = 110000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 110000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 110000
$ws.Range("N137").Value = -120200

$ws.Range("H138").Value = 20390
$ws.Range("I138").Value = 20390
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 20390
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -15250

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0

$ws.Range("H140").Value = 40000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 40000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

